# bioSample_3275.xlsx — "updated fastq files and related metadata"
#
# 1) treatment column (G) changes from "37C.CO2" to "DMEM.37C.CO2" for every
#    existing data row (rows 2-37).
# 2) four new samples are appended as rows 38-41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) update the treatment column for all existing rows ------------------
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 7).Value = "DMEM.37C.CO2"
}

# --- 2) append the four new sample rows -------------------------------------
# Column A (harvestDate) holds text that looks like a date (MM.DD.YY); a
# leading apostrophe keeps Excel from reinterpreting it as a real date, and
# resetting the style afterwards clears the quote-prefix formatting flag so
# the cell ends up as plain text, matching the rest of the column.

function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 38
Set-TextCell $ws.Cells.Item(38, 1) "08.09.18"
$ws.Cells.Item(38, 2).Value = "H.BROWN"
$ws.Cells.Item(38, 3).Value = 38
$ws.Cells.Item(38, 4).Value = "90minuteinduction"
$ws.Cells.Item(38, 5).Value = "TDY2205"
$ws.Cells.Item(38, 6).Value = "CNAG_06871"
$ws.Cells.Item(38, 7).Value = "DMEM.37C.CO2"
$ws.Cells.Item(38, 8).Value = 90
$ws.Cells.Item(38, 9).Value = 5

# Row 39
Set-TextCell $ws.Cells.Item(39, 1) "10.15.18"
$ws.Cells.Item(39, 2).Value = "H.BROWN"
$ws.Cells.Item(39, 3).Value = 39
$ws.Cells.Item(39, 4).Value = "90minuteinduction"
$ws.Cells.Item(39, 5).Value = "TDY1452"
$ws.Cells.Item(39, 6).Value = "CNAG_02566"
$ws.Cells.Item(39, 7).Value = "DMEM.37C.CO2"
$ws.Cells.Item(39, 8).Value = 90
$ws.Cells.Item(39, 9).Value = 13

# Row 40
Set-TextCell $ws.Cells.Item(40, 1) "10.30.18"
$ws.Cells.Item(40, 2).Value = "H.BROWN"
$ws.Cells.Item(40, 3).Value = 40
$ws.Cells.Item(40, 4).Value = "90minuteinduction"
$ws.Cells.Item(40, 5).Value = "TDY1118"
$ws.Cells.Item(40, 6).Value = "CNAG_05222"
$ws.Cells.Item(40, 7).Value = "DMEM.37C.CO2"
$ws.Cells.Item(40, 8).Value = 90
$ws.Cells.Item(40, 9).Value = 1

# Row 41
Set-TextCell $ws.Cells.Item(41, 1) "10.15.18"
$ws.Cells.Item(41, 2).Value = "H.BROWN"
$ws.Cells.Item(41, 3).Value = 41
$ws.Cells.Item(41, 4).Value = "90minuteinduction"
$ws.Cells.Item(41, 5).Value = "TDY1174"
$ws.Cells.Item(41, 6).Value = "CNAG_00871"
$ws.Cells.Item(41, 7).Value = "DMEM.37C.CO2"
$ws.Cells.Item(41, 8).Value = 90
$ws.Cells.Item(41, 9).Value = 10

# --- 3) selection now spans the whole A:K column range ----------------------
$ws.Range("A1:K1048576").Select()

# --- 4) window view was reset/maximized -------------------------------------
$win = $excel.Windows.Item(1)
$win.Left = 0
$win.Top = 0
$win.Width = 33600
$win.Height = 21000
